$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "science" (token value 58) as a new row right before the old
# row 20 ("{" / 15). Rows.Insert shifts rows 20-45 down to 21-46, and
# (per observed engine behaviour) the cells pushed down keep their
# original per-row height/format, while the new row inherits style s="1"
# plus full A:O population, matching rows 1-19 immediately above it.
$ws.Rows(20).Insert()
$ws.Range("B20").Value = "science"
$ws.Range("C20").Value = 58

# --- Insert ":" (token value 60) as a new row right before the (now
# shifted) ";" row, which after the previous insert sits at row 35.
$ws.Rows(35).Insert()
$ws.Range("B35").Value = ":"
$ws.Range("C35").Value = 60

# --- Append the two brand new tokens "new" (57) and "SuchWow" (59) at
# the bottom of the table, right after the current last row (47, the
# shifted '"' row). Use an Insert-based approach so the new rows inherit
# the minimal 2-column (B:C) layout used by the tail of the table, then
# pin the row height to match its neighbours.
$ws.Rows("48:49").Insert()
$ws.Range("B48").Value = "new"
$ws.Range("C48").Value = 57
$ws.Range("B49").Value = "SuchWow"
$ws.Range("C49").Value = 59
$ws.Rows("48:49").RowHeight = 18

# --- Restore the view state Excel would naturally leave after scrolling
# to / selecting near the newly appended rows.
$ws.Range("A32").Select()
$excel.ActiveWindow.ScrollRow = 32
$ws.Range("C52").Select()
